$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 4
$scratchRow = 100

function Swap-Rows($ws, $r1, $r2, $lastCol, $scratchRow) {
    $rowA = $ws.Range($ws.Cells.Item($r1, 1), $ws.Cells.Item($r1, $lastCol))
    $rowB = $ws.Range($ws.Cells.Item($r2, 1), $ws.Cells.Item($r2, $lastCol))
    $scratch = $ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow, $lastCol))

    # Use Range.Copy (not .Value/.Formula) so cell types/text are preserved
    # verbatim -- numeric-looking IDs like "645740.1" stay text instead of
    # being re-parsed into numbers, and no quotePrefix style gets added.
    $scratch.ClearContents()
    $rowA.Copy($scratch)

    $rowA.ClearContents()
    $rowB.Copy($rowA)

    $rowB.ClearContents()
    $scratch.Copy($rowB)

    $scratch.ClearContents()

    # Copy() of a blank source cell leaves a bare empty <c/> placeholder
    # behind instead of truly clearing the cell -- sweep once more so any
    # cell that should end up empty is actually removed from the sheet XML.
    foreach ($rng in @($rowA, $rowB)) {
        foreach ($cell in $rng.Cells) {
            if ($cell.Formula -eq "") {
                $cell.ClearContents()
            }
        }
    }
}

# Swap NOS2P2 (row 2) and NOS2P1 (row 3)
Swap-Rows $ws 2 3 $lastCol $scratchRow

# Swap SDR42E2 (row 12) and AKR1C8P (row 13) - fixes the Alias ("AKR1CL1")
# and Name so they line up with the correct Symbol/model instead of
# referencing the non-existing one.
Swap-Rows $ws 12 13 $lastCol $scratchRow

# Swap COX2 (row 16) and COX1 (row 17)
Swap-Rows $ws 16 17 $lastCol $scratchRow

# Remove the scratch row entirely so the sheet's used range/dimension goes
# back to its original extent.
$ws.Rows.Item($scratchRow).Delete()
